$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The "antiguedad / rehabilitacion / calidad_alta / estado_conservacion"
# question block (columns H:K) is retired and replaced by a new
# "calidad_alta / calefaccion / creci / renta" block.
#
# Writes below are ordered so that brand-new shared strings land in the
# workbook's shared-string table in the same order they appear in the
# updated template: calefaccion, Sí, creci, renta.
# -----------------------------------------------------------------------

# Header row
$ws.Range("H1").Value = "calidad_alta"
$ws.Range("I1").Value = "calefaccion"
$ws.Range("I2").Value = "Sí"
$ws.Range("J1").Value = "creci"
$ws.Range("K1").Value = "renta"

# Row 2
$ws.Range("H2").Value = $false
$ws.Range("K2").Value = "Sí"

# Row 3
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = "Sí"
$ws.Range("K3").Value = "No"

# Row 4
$ws.Range("H4").Value = $false
$ws.Range("I4").Value = "Sí"
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = "No"

# Row 5
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = "No"
$ws.Range("J5").Value = "Sí"
$ws.Range("K5").Value = $false

# Re-fit the columns that just had their headers/values replaced.
$ws.Columns("H:K").AutoFit()

# -----------------------------------------------------------------------
# View state: the sheet is now viewed at 100% zoom with a different cell
# selected.
# -----------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("I12").Select() | Out-Null
